$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell references whose new values are numeric-looking strings that must
# stay as Text (matching the source data, which stores prices as text).
# Force a Text number format before assigning, then restore the default
# "Normal" style so no stray formatting is introduced.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '68.418.10'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '2.648.97'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue $ws.Range("D5") '597.36'
$ws.Range("E5").Value = '  -0.17%  '
Set-TextValue $ws.Range("D6") '158.90'
$ws.Range("E6").Value = '  +2.67%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("D9").Value = '2.648.23'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("E10").Value = '  -2.06%  '
$ws.Range("E11").Value = '  -1.01%  '
Set-TextValue $ws.Range("D12") '5.28'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("E13").Value = '  -0.85%  '
Set-TextValue $ws.Range("D14") '28.02'
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").Value = '3.132.91'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("E16").Value = '  -2.97%  '
$ws.Range("D17").Value = '68.301.11'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '2.644.08'
$ws.Range("E18").Value = '  -0.27%  '
Set-TextValue $ws.Range("D19") '11.55'
$ws.Range("E19").Value = '  +1.44%  '
Set-TextValue $ws.Range("D20") '363.91'
$ws.Range("E20").Value = '  -0.23%  '
Set-TextValue $ws.Range("D21") '7.47'
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("E23").Value = '  -2.10%  '
Set-TextValue $ws.Range("D24") '2.08'
$ws.Range("E24").Value = '  +1.12%  '
Set-TextValue $ws.Range("D25") '74.72'
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("E26").Value = '  +0.02%  '
Set-TextValue $ws.Range("D27") '9.89'
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("D28").Value = '2.783.90'
$ws.Range("E28").Value = '  +0.18%  '
Set-TextValue $ws.Range("D29") '0.0000103'
$ws.Range("E29").Value = '  -3.13%  '
Set-TextValue $ws.Range("D30") '1.00'
$ws.Range("E30").Value = '  +0.37%  '
Set-TextValue $ws.Range("D31") '565.10'
$ws.Range("E31").Value = '  +0.29%  '
Set-TextValue $ws.Range("D32") '8.07'
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("E33").Value = '  -0.49%  '
$ws.Range("E34").Value = '  +0.36%  '
Set-TextValue $ws.Range("D35") '1.65'
$ws.Range("E35").Value = '  +4.39%  '
$ws.Range("E36").Value = '  -1.79%  '
$ws.Range("E37").Value = '  +0.00%  '
Set-TextValue $ws.Range("D38") '160.69'
$ws.Range("E38").Value = '  -0.34%  '
Set-TextValue $ws.Range("D39") '19.65'
$ws.Range("E39").Value = '  +1.66%  '
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("E41").Value = '  -0.88%  '
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("E44").Value = '  -5.74%  '
Set-TextValue $ws.Range("D46") '158.23'
$ws.Range("E46").Value = '  +1.20%  '
$ws.Range("E47").Value = '  +1.71%  '
Set-TextValue $ws.Range("D48") '21.98'
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("E51").Value = '  +1.85%  '
